$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 72.59633679618487
$ws.Range("C2").Value = 75.4195437681842
$ws.Range("D2").Value = 65.68624930906326
$ws.Range("E2").Value = 77.98416199942649
$ws.Range("B3").Value = 95.43563379966095
$ws.Range("C3").Value = 94.57003139809234
$ws.Range("D3").Value = 95.41026809330305
$ws.Range("E3").Value = 95.18253402370901
$ws.Range("B4").Value = 99.38626481381891
$ws.Range("C4").Value = 99.2978203857076
$ws.Range("D4").Value = 99.37140167670691
$ws.Range("E4").Value = 99.39704287525669
$ws.Range("B5").Value = 98.90471815478253
$ws.Range("C5").Value = 98.91138032066465
$ws.Range("D5").Value = 98.8905065767255
$ws.Range("E5").Value = 98.8626901422292
$ws.Range("B6").Value = 98.49158790473874
$ws.Range("C6").Value = 98.41221242544741
$ws.Range("D6").Value = 98.4250015907633
$ws.Range("E6").Value = 98.3711814400055
$ws.Range("B7").Value = 97.9066693623124
$ws.Range("C7").Value = 97.93034118333546
$ws.Range("D7").Value = 97.9708024524988
$ws.Range("E7").Value = 97.93196603237055
$ws.Range("B8").Value = 97.44840259030759
$ws.Range("C8").Value = 97.40936359917099
$ws.Range("D8").Value = 97.46402796050199
$ws.Range("E8").Value = 97.41500342746428
$ws.Range("B9").Value = 96.06184477345799
$ws.Range("C9").Value = 96.0611824867761
$ws.Range("D9").Value = 96.04396408512935
$ws.Range("E9").Value = 96.09274653964862
